# C5-PowerPoint.pptx edit — Fri, Jul 24, 2020 5:05:16 PM
#
# 1) Slide 6's table switches to a different built-in PowerPoint table
#    style (tableStyleId {209203A3-...} -> {E363592B-...}).
# 2) The deck's main theme (theme1.xml, used by the slide master) is
#    recoloured from the "Integral" palette to the standard Office
#    palette (the companion theme2.xml, used only by the notes master,
#    already carries the complementary palette).

$p = $ppt.ActivePresentation

# Helper: convert an "RRGGBB" hex string (as it appears in DrawingML
# <a:srgbClr val="..."/>) into the BGR-packed long that the PowerPoint
# RGB color properties expect over COM.
function ToCOMColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

# --- 1) Table style on slide 6 -------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{E363592B-F288-4C24-9FC5-DAEA573DDB49}", $true)
    }
}

# --- 2) Recolour the theme used by the slide master to the Office palette -
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# Order matches MsoThemeColorSchemeIndex: dk1, lt1, dk2, lt2, accent1-6,
# hyperlink, followed hyperlink.
$officePalette = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hyperlink
    "954F72"  # followed hyperlink
)

for ($i = 1; $i -le $officePalette.Count; $i++) {
    $colorScheme.Item($i).RGB = ToCOMColor $officePalette[$i - 1]
}
